$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2 and 3: column C time ranges change (column B text unchanged)
$ws.Cells.Item(2, 3).Value = '10:55-11:0'
$ws.Cells.Item(3, 3).Value = '11:0-11:5'

# Add new rows 4-15 with encrypted message (col B) and time range (col C)
$ws.Cells.Item(4, 2).Value = 'U2FsdGVkX1+EqRz5mg36pN1B7kT9TR74WNMW/ISHW6cmSoQSnp5jVxizWU7v40DIEcFLRpkmAsuMakXso2huErv5jipCq63FnUFQEFKQXScxoglLuNsCEAroagHty+NcVBAK06ajVBZD9FZYUo1INlDEzJ5CdHvJHxlfTVzdd+cGGAVluDv+kdxOtAKpnGb5oGAYnJfu2OtdcKykwKqfhkELCcAxJNUuDgoOR9B4OGEM62hFHbG6Yw62/AspDOUvsPz7gMc2JOztQnFmoWfvKWs5CIW2fF9ZhDt/eTrtboRnByBZE15uaFyGLnxvpu+oUaipP/9ecNG0NV7nFt/BDLEClRzUmKMLCaUgv2qXmIuLXh4KyJvl/DSUSx2fIAvEdXwI619kb8AMwAUxMY6nTq6Xxl4cDSOQxrGOzUe/i7Q='
$ws.Cells.Item(4, 3).Value = '11:5-11:10'
$ws.Cells.Item(5, 2).Value = 'U2FsdGVkX18D85up224eLUPZ7GpNJAtvdfgU27dFmQInnqsUmJZEyVW6bReRwHft383DfZIRXHbYx50rmTipM8dTB2HjIk7aD/UkTIXXZ9RrXmBr4oSSN3Q2l2wkGmmrUuPC9tWFLDSUMhWesI2dtAxd/G82Qo1Pxunprz0gRo3pivMSOYpw5RNIDP2K8Z9itAsEEB3icV+HzUhYDeUM8Ks91PnJE9NOrZ2t74NLzStAS/mh38dzgkX30JK3hZArbKISggzxsW5stnK/0h+tqIEokEWkDZVMCAaXeoL38fh26+ToJH8AoDngD4LODjTTs0T6C6VRPa/BbZhc55yZpFUpK5ycV7P2Ci57H2Ac3dU='
$ws.Cells.Item(5, 3).Value = '11:10-11:15'
$ws.Cells.Item(6, 2).Value = 'U2FsdGVkX19FxGzlsvUgo+buXRdS1RLU7kJMDzGAdMUQkckbrgbWv7jz+QT/6N9ZdX16Z+EenDeUyTV2nhC3w3vZTq0kRWc2vK/u3kAH70oABJzDDxciUiptkf/q4t8rRWthXtzJPMsFx4NnaAtQlQPRmwtHlIcT9uZm3x1QhtkZV7+VWOS1nvKUX4VSXAZ0JVJNZ59f1+K83wFQQLRSYurvXEba0gJ7LAiBSRhs0pi/nslCJO9gSatF4jCGpDawRVGCqw3bjCwz0rDZKuvr63M8k/KubEo1BXYd1nWbe+1mhOxEod948hhOLz8XveQnf3L5QhOcArrg0RMqI1PxfShEz7+ka0lkv7MFm0Qj63M8s12l9v5/b9r1j9Z3V52ZJ7BaOCAA4kY3yju27j/ie1M2JP5DVMMy8grcUza7M/DIGb1zoSZ8zjosBZArClSy'
$ws.Cells.Item(6, 3).Value = '11:15-11:20'
$ws.Cells.Item(7, 2).Value = 'U2FsdGVkX18dvSQioA8bbZXBmKNPWCxwHchn7izQU9COIhFMq1KX36lQNjZiYrQD//S8EmsoDReRaH5yElamMDA2wrD0xDVMM9G5Z9redIxAtPyF0Ag0NV2Z3sKCmzeh58n7mLUu6y10ygW+qt+uSwkbjaOG9780adK1azYqeGdIiLDyOQqiuo24zZwAa6iisZzP0Lykag76rwMfjhU3pav3rwswzF+zgm0uA83ezdnlD0H7eNxENQWFu+dWQSL+ERdiVkAMQ8tVCF4cxvCf6Pe35oFj77EIB30Y8y3wFVYy+GvAJOjhaIuo/FdNmwtMxm5u2ypO09dhXcHq4RQa1yKzVKqqX1sglTDW8V7qON2PN1pWynvzeBlMHWfNuDPNPT1ejb+nnIzY4yNqEOWEVRywe0vS1ICO4ft/3yEJF+s='
$ws.Cells.Item(7, 3).Value = '11:20-11:25'
$ws.Cells.Item(8, 2).Value = 'U2FsdGVkX19eQZGX7pj68kIum1I0o+zBh4WPYNf8xq3siizN52bKew5l+nywF74gq5cAzJ9tkn51FYLRRTGi66vuDCqq1eOjCI7bum6vRtm8s4UcObcZXRpujZx2sDq/KIEiSRguShI3VOEhk+WmZqLqwG7tOqBNJZXz6lW+/oA7NRpChVvxhWFNfTgwOwZmPPjmhgZqTlVz5GdAHI0mYifDY6lyvD1aND2PKn5ZSmxWhPwzMXmVOEojm8tnNWreCmsh1aqrKA7+FCppwudh0gDp1iNrNhDG20sYZz6bymF7PoqDhOQ8vnl3kCBUxGfhbg5gQANqNprC09QxsEx/wdi+M7D5kzEuvLQX8CxIkzAqbkPVG3W8v/uR0QFFtJfz3Z7N/MqV3XCt7dInXL1k9ilnxmifWyUz8k0RHCCfoOnpNyQdcYjWKutaHCWi2gyZyd/NFmZWOS1nl40omEvjMg=='
$ws.Cells.Item(8, 3).Value = '11:25-11:30'
$ws.Cells.Item(9, 2).Value = 'U2FsdGVkX19X3MgjX8KTTcrESJoNSwdMSlfC20xYZa42RtbdCHOOY/xohgjv2dtNrUbhhBvHQqEazTNLBhJNlNQSGCVC1zcIe3TKFuDcc1o+OIirIDfGvYpS2rRl1SWPIP6SQSHhc+cy4+VlTbmrIE6jOK84QRtI8f0bEzx95hLWlkQgARy17z3SFtaUToyAGkBJtC85swugMR/lvLNwgZFEk0TCuCUsMs3rjpPdQUhE9Wj9K9sOFeVXB3j+0lbO71XA+P5P6w6bUfa3Be/ov/Prp4ngaazvuG/gO4caCV5/IqJzVLxzpotSgX7wFrZNBEmLm8kdWm1DfKh4muju8JQFhMeu7WE9l4aJ8nSXsUY='
$ws.Cells.Item(9, 3).Value = '11:30-11:35'
$ws.Cells.Item(10, 2).Value = 'U2FsdGVkX19vF2FV0aQYAgE7L+9745pjetW0kMg3EzbL+IVHC27R1HpZ2zFDxVUMAKOuFGVJ8OXT0+XXOjrkMDv5u/Y6f/8wHxctZz1lW6H0sEgH/VPjtQ6meyZmuinIS+eq5TVZkTk9QUjQ23dxllHNba0kgTSStYY+4wnKeFafT2M17BXi45d6HU0ig4YmclM7dHDnxVHK/Z+bym7qCNQ2Pzz4mnNjMA0aJEhV3cOR0sbsvjGZ8XDiWyyqQvKROqOGp5tPk0n5IrBH4yWiH8/Qw8hcbIpvwoa5yO42kX0ILZ0YkUwIWvWoI2u8Cma6qE2I3eXaE6BlF81NxHQ6E/7woHHtVPQkFBSm9Dj5K2A='
$ws.Cells.Item(10, 3).Value = '11:35-11:40'
$ws.Cells.Item(11, 2).Value = 'U2FsdGVkX1+v2CH9e8k+TdlapGKj+O/wO2ylWEfgkOmZHO2fFAxLjna5VA50nTtNUifZMfuMihNm+KW3Vp/DPYSM3OvJljKIcWILvxEnt4Uq9FYzHSk76l2SFXn9GWG80snjjvXZtQwX7Re79C8/X3JHFyUA0p1NsGXhSVwSPANVaoWqPGIO9VlrRUWyP36N5dkRXy17hQ8N0CeL+Lp3qPhwhazfOaOu7DewaLaS9LYzTnnJijJeFq06VMtnOcepIWm14H1wPOT5fw5fSP87w+CZrwCs6fDB7w1sgO3nlgWCtrswltYg+rdH3LWIu5aCfIV9VLqPTcc45MIOQTcLf3fXeK6ZZmMsnyCxFWY5U7I/9N8+B0qyh8nvQzHMjpm7Vq8kAIuEOVIgUYWtoYmxR6FAIKD8iT0zrF8itzvOF6Vte8AezZOradPp1nABFHOCRNl2iloWRXVm5GLWLpNyfUwfYZTZ3NiOzMIYFJjUKeg='
$ws.Cells.Item(11, 3).Value = '11:40-11:45'
$ws.Cells.Item(12, 2).Value = 'U2FsdGVkX19UVYawAGmNToFy56xASB15xGyt4+Ehl44thvWSNfdqyFTGhHsz0aQhR3Z0wjCumkE3FGV1dwwe8hbqHxLnbJ2igAkhZCD9Pu4YQJ+sai/c3mHS2RWqr5++3bB//2RD4uMOvNv9ill6Cd40gEBzwBz4E/eKlKt8xe3aXUZsltM7D/+qoTqp7xMbvr/NPl+bwim4LUJu4umutGs0Xt255rDhGckL2qmJo4qM/S7bNDfrH/W8E6UcZ0FyLz+vuJG0YwedHloP/+lmLQUWwjGhQK5xl8VQ33hzgUzBPdb99wnlfQKLYJ4jYNKj7XPtA2qcR7sn58PT0vsjs8OjzDbsgvkZDStBFcK9n/NzYHRgLIrcudodqIT1IYRTRhJ5mADZppYERWe7OBRWcQ=='
$ws.Cells.Item(12, 3).Value = '11:45-11:50'
$ws.Cells.Item(13, 2).Value = 'U2FsdGVkX18VQxeyU6dxfSKRf/WHcLalX6DnH62kdRPfx0KBjjSzopcbJ7fc1f/msjb28WJX6nj+AiAG9v0cYJ7aKd/zHjrfB3SRnKEpaXWs+Hn/kix4N14UTeAkDo5FEUsyCu6XqNDBHlU+gDLWTzYX3bcY3DaZG9UVXoba4Q3jEJWJDedL42Eim4HFOnE5gUIbGwyz7+GCAE4njctTuiY3FpL0TA4WcI/fYTPowpLD7ZLtCqPXgs4Pvo34umBnEaTQucajC5jvvfJIa84G0fbVD1MtjvknJm2T+m5OrtaTg7EPOHoClgRNCefYSp+LpLWZ3d6weOr3NJ6+igv0+C1Zla52kjlDLwlfyCE/Bo+elQP6xJLjZRlZs8/cgfqpECd3TM5B+IXUYJEnD1bD5/Q4pvXiSMI+qbpakLMaPkYdH7w/SlqQueSWpOPMh7D3mKMoDfmo8q2GjLOYCjVSPg=='
$ws.Cells.Item(13, 3).Value = '11:50-11:55'
$ws.Cells.Item(14, 2).Value = 'U2FsdGVkX1+8KKxqt5fT1eBZf8nUskNJ7mFCE1QIKREeppPwhUh1Sa0aWqiglHuc'
$ws.Cells.Item(14, 3).Value = '11:55-12:0'
$ws.Cells.Item(15, 2).Value = 'U2FsdGVkX1/uSd6oYz/AMEPgHPLjc2amon7N0LXziHgeKMagBUQlnbder3K8r2geGkVIi6HvkoK7S2QaPfRdZqXKEQGAI7ohJX1P955fY9fgkD5keunwkq86pzZIuMJJ/TjEydLCnSxP2E6Sn+G+7gHG/7Clo+Hn1VnGCnAYh8Etw/pcxfzyE5c5YJpdQtpIWHxowWDP1+y7WGANO+H1qQ=='
$ws.Cells.Item(15, 3).Value = '12:0-12:5'

# Update selection to match the saved view state
$ws.Range("B18").Select()
